$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for data rows 2-7 (columns D,K,L,M,N,O,P,R,S).
# All other columns (A,B,C,E,F,G,H,I,J,Q,T) are unchanged.
$rows = @(
    @{ Row = 2;  D = 44305; K = "Mankaki"; L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 },
    @{ Row = 3;  D = 44699; K = "Mankaki"; L = "Primera"; M = 250; N = 29000; O = 30000; P = 29500; R = "Región de O'Higgins"; S = 1639 },
    @{ Row = 4;  D = 44313; K = "Mankaki"; L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; R = "Región de O'Higgins"; S = 1194 },
    @{ Row = 5;  D = 44301; K = "Hachiya"; L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; R = "Región de O'Higgins"; S = 1139 },
    @{ Row = 6;  D = 44355; K = "Mankaki"; L = "Segunda"; M = 270; N = 20000; O = 21000; P = 20500; R = "Región Metropolitana"; S = 1139 },
    @{ Row = 7;  D = 44342; K = "Mankaki"; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
}
